$d = $word.ActiveDocument

# Locate the "search feature" bullet paragraph and recolor its text (and
# paragraph mark) dark green (RRGGBB 008000), matching wdColorDarkGreen.
$wdColorDarkGreen = 0x008000
$target = "Írnod kell egy keresési funkciót az oldalhoz."

$matchCount = 0
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*$target*") {
        $r.Font.Color = $wdColorDarkGreen
        $matchCount = $matchCount + 1
    }
}

Write-Output "Recolored $matchCount paragraph(s)."
